$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 57884

# Row 3
$ws.Range("B3").Value = 57884

# Row 4
$ws.Range("B4").Value = 57884

# Row 5
$ws.Range("B5").Value = 57884

# Row 6
$ws.Range("B6").Value = 80252

# Row 7
$ws.Range("A7").Value = 131009283
$ws.Range("B7").Value = 79243
$ws.Range("E7").Value = 6425
$ws.Range("F7").Value = 'Garnlav'
$ws.Range("G7").Value = 'Alectoria sarmentosa'
$ws.Range("H7").Value = '(Ach.) Ach.'
$ws.Range("M7").Value = $null
$ws.Range("Q7").Value = 589968
$ws.Range("R7").Value = 6911120
$ws.Range("Z7").Value = '11:03'
$ws.Range("AB7").Value = '11:03'
$ws.Range("AC7").Value = $null

# Row 8
$ws.Range("B8").Value = 57884

# Row 9
$ws.Range("A9").Value = 131009297
$ws.Range("B9").Value = 57884
$ws.Range("Q9").Value = 589752
$ws.Range("R9").Value = 6911214
$ws.Range("Z9").Value = '10:18'
$ws.Range("AB9").Value = '10:18'

# Row 10
$ws.Range("A10").Value = 131009270
$ws.Range("B10").Value = 57884
$ws.Range("E10").Value = 100109
$ws.Range("F10").Value = 'Tretåig hackspett'
$ws.Range("G10").Value = 'Picoides tridactylus'
$ws.Range("H10").Value = '(Linnaeus, 1758)'
$ws.Range("M10").Value = 'färska spår'
$ws.Range("Q10").Value = 589700
$ws.Range("R10").Value = 6911274
$ws.Range("Z10").Value = '12:22'
$ws.Range("AB10").Value = '12:22'
$ws.Range("AC10").Value = 'färska ringhack på tall'

# Row 11
$ws.Range("B11").Value = 57884

# Row 12
$ws.Range("B12").Value = 80221

# Row 13
$ws.Range("B13").Value = 57884

# Row 14
$ws.Range("B14").Value = 57884

# Row 15
$ws.Range("B15").Value = 80252

# Row 16
$ws.Range("B16").Value = 57884

# Row 17
$ws.Range("A17").Value = 131009301
$ws.Range("B17").Value = 57884
$ws.Range("Q17").Value = 589700
$ws.Range("R17").Value = 6911168
$ws.Range("Z17").Value = '09:59'
$ws.Range("AB17").Value = '09:59'

# Row 18
$ws.Range("A18").Value = 131009282
$ws.Range("B18").Value = 57884
$ws.Range("Q18").Value = 589986
$ws.Range("R18").Value = 6911103
$ws.Range("Z18").Value = '11:12'
$ws.Range("AB18").Value = '11:12'

# Row 19
$ws.Range("B19").Value = 57884

# Row 20
$ws.Range("B20").Value = 57884

# Row 21
$ws.Range("B21").Value = 57884

# Row 22
$ws.Range("A22").Value = 131009308
$ws.Range("B22").Value = 79243
$ws.Range("E22").Value = 6425
$ws.Range("F22").Value = 'Garnlav'
$ws.Range("G22").Value = 'Alectoria sarmentosa'
$ws.Range("H22").Value = '(Ach.) Ach.'
$ws.Range("I22").Value = $null
$ws.Range("M22").Value = $null
$ws.Range("Q22").Value = 589686
$ws.Range("R22").Value = 6911077
$ws.Range("Z22").Value = '09:17'
$ws.Range("AB22").Value = '09:17'

# Row 23
$ws.Range("B23").Value = 57884

# Row 24
$ws.Range("B24").Value = 57884

# Row 25
$ws.Range("A25").Value = 131009304
$ws.Range("B25").Value = 58043
$ws.Range("E25").Value = 103021
$ws.Range("F25").Value = 'Talltita'
$ws.Range("G25").Value = 'Poecile montanus'
$ws.Range("H25").Value = '(Conrad von Baldenstein, 1827)'
$ws.Range("I25").Value = '1'
$ws.Range("M25").Value = 'lockläte, övriga läten'
$ws.Range("Q25").Value = 589753
$ws.Range("R25").Value = 6911167
$ws.Range("Z25").Value = '09:45'
$ws.Range("AB25").Value = '09:45'

# Row 26
$ws.Range("B26").Value = 57884

# Row 27
$ws.Range("B27").Value = 57884

# Row 28
$ws.Range("B28").Value = 57884

# Row 29
$ws.Range("B29").Value = 57884

# Row 30
$ws.Range("A30").Value = 131009285
$ws.Range("B30").Value = 80252
$ws.Range("D30").Value = 'LC'
$ws.Range("E30").Value = 6456
$ws.Range("F30").Value = 'Skinnlav'
$ws.Range("G30").Value = 'Leptogium saturninum'
$ws.Range("H30").Value = '(Dicks.) Nyl.'
$ws.Range("M30").Value = $null
$ws.Range("Q30").Value = 589879
$ws.Range("R30").Value = 6911153
$ws.Range("Z30").Value = '10:56'
$ws.Range("AB30").Value = '10:56'
$ws.Range("AC30").Value = $null

# Row 31
$ws.Range("A31").Value = 131009302
$ws.Range("B31").Value = 57884
$ws.Range("Q31").Value = 589716
$ws.Range("R31").Value = 6911140
$ws.Range("Z31").Value = '09:54'
$ws.Range("AB31").Value = '09:54'
$ws.Range("AC31").Value = 'färska ringhack på tall'

# Row 32
$ws.Range("A32").Value = 131009306
$ws.Range("B32").Value = 57884
$ws.Range("D32").Value = 'NT'
$ws.Range("E32").Value = 100109
$ws.Range("F32").Value = 'Tretåig hackspett'
$ws.Range("G32").Value = 'Picoides tridactylus'
$ws.Range("H32").Value = '(Linnaeus, 1758)'
$ws.Range("M32").Value = 'färska spår'
$ws.Range("Q32").Value = 589791
$ws.Range("R32").Value = 6911148
$ws.Range("Z32").Value = '09:38'
$ws.Range("AB32").Value = '09:38'
$ws.Range("AC32").Value = 'färska och äldre ringhack på tall'

# Row 33
$ws.Range("B33").Value = 57884

# Row 34
$ws.Range("B34").Value = 57884

# Row 35
$ws.Range("B35").Value = 57884

# Row 36
$ws.Range("B36").Value = 57884

# Row 37
$ws.Range("B37").Value = 79243

# Row 38
$ws.Range("B38").Value = 57884

# Row 39
$ws.Range("B39").Value = 57884

# Row 40
$ws.Range("A40").Value = 131009291
$ws.Range("B40").Value = 80252
$ws.Range("D40").Value = 'LC'
$ws.Range("E40").Value = 6456
$ws.Range("F40").Value = 'Skinnlav'
$ws.Range("G40").Value = 'Leptogium saturninum'
$ws.Range("H40").Value = '(Dicks.) Nyl.'
$ws.Range("M40").Value = $null
$ws.Range("Q40").Value = 589791
$ws.Range("R40").Value = 6911200
$ws.Range("Z40").Value = '10:37'
$ws.Range("AB40").Value = '10:37'
$ws.Range("AC40").Value = $null

# Row 41
$ws.Range("A41").Value = 131009275
$ws.Range("B41").Value = 57884
$ws.Range("D41").Value = 'NT'
$ws.Range("E41").Value = 100109
$ws.Range("F41").Value = 'Tretåig hackspett'
$ws.Range("G41").Value = 'Picoides tridactylus'
$ws.Range("H41").Value = '(Linnaeus, 1758)'
$ws.Range("M41").Value = 'färska spår'
$ws.Range("Q41").Value = 589844
$ws.Range("R41").Value = 6911365
$ws.Range("Z41").Value = '11:53'
$ws.Range("AB41").Value = '11:53'
$ws.Range("AC41").Value = 'färska ringhack på tall'
